# Update the public EPEX Spot prices workbook:
#  - "Prix Spot" sheet: insert a new date column (07-nov) before the
#    existing 01-oct. column (column DL), filling the new column with
#    "-" placeholders for every hour row (the day's data isn't in yet).
#  - "Gaz" / "CO2" sheets: append the next day's price row (2025-11-05).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Prix Spot" - insert the 07-nov column at DL (column 116),
# shifting the former DL:EP ("01-oct." .. "31-oct.") block one column
# right to DM:EQ.
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Columns.Item(116).Insert()

$wsPrix.Range("DL1").Value = "07-nov"
$wsPrix.Range("DL2:DL25").Value = "-"

# ---------------------------------------------------------------------
# Sheet 2: "Gaz" - append row 144 with the 2025-11-05 price.
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$cellGaz = $wsGaz.Range("A144")
$cellGaz.NumberFormat = "@"
$cellGaz.Value = "2025-11-05"
$wsGaz.Range("B144").Value = 30.425

# ---------------------------------------------------------------------
# Sheet 3: "CO2" - append row 144 with the 2025-11-05 price.
# ---------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$cellCo2 = $wsCo2.Range("A144")
$cellCo2.NumberFormat = "@"
$cellCo2.Value = "2025-11-05"
$wsCo2.Range("B144").Value = 81.18000000000001
